# "update cloud related figures"
#
# 1. dt!C8 sensor label changes from "MYD" to "MYD10" (a new, more specific
#    sensor label). This retires the shared string "MYD" (freeing it from the
#    table) and introduces "MYD10".
# 2. The selection on the "dt" sheet moves to E12.
# 3. A brand-new "Sheet1" worksheet is appended after "MYD" holding a small
#    AWS cloud-related comparison table (stations x years) with a red-yellow
#    -green color-scale conditional format, and becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: tweak the "dt" sheet -------------------------------------------
$dt = $wb.Worksheets.Item("dt")
$dt.Range("C8").Value = "MYD10"
$dt.Range("E12").Select()

# --- 3: build the new "Sheet1" worksheet, inserted after the last sheet ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)

# header row
$ws.Range("A1").Value = "AWS"
$ws.Range("B1").Value = 2020
$ws.Range("C1").Value = 2021
$ws.Range("D1").Value = 2022

# data rows -- written in this particular order so newly-introduced shared
# strings land in the same order as the authored workbook
$ws.Range("A4").Value = "kpc_l"
$ws.Range("B4").Value = 0.64
$ws.Range("C4").Value = 0.88
$ws.Range("D4").Value = 0.84
$ws.Range("A5").Value = "nuk_k"
$ws.Range("B5").Value = 0.42
$ws.Range("C5").Value = 0.43
$ws.Range("D5").Value = 0.44
$ws.Range("A6").Value = "nuk_l"
$ws.Range("B6").Value = 0.08
$ws.Range("C6").Value = 0.36
$ws.Range("D6").Value = 0.47
$ws.Range("A7").Value = "qas_l"
$ws.Range("B7").Value = 0.14
$ws.Range("C7").Value = 0.41
$ws.Range("D7").Value = 0.15
$ws.Range("A8").Value = "qas_m"
$ws.Range("B8").Value = 0.92
$ws.Range("C8").Value = 0.58
$ws.Range("D8").Value = 0.89
$ws.Range("A9").Value = "qas_u"
$ws.Range("B9").Value = 0.67
$ws.Range("C9").Value = 0.7
$ws.Range("D9").Value = 0.67
$ws.Range("A10").Value = "tas_l"
$ws.Range("B10").Value = -0.09
$ws.Range("C10").Value = 0.67
$ws.Range("D10").Value = 0.59
$ws.Range("A11").Value = "thu_l"
$ws.Range("B11").Value = 0.74
$ws.Range("C11").Value = 0.4
$ws.Range("D11").Value = 0.81
$ws.Range("A12").Value = "thu_u"
$ws.Range("B12").Value = 0.63
$ws.Range("C12").Value = 0.15
$ws.Range("D12").Value = 0.65
$ws.Range("A13").Value = "thu_u2"
$ws.Range("B13").Value = 0.4
$ws.Range("C13").Value = 0.31
$ws.Range("D13").Value = 0.05
$ws.Range("A14").Value = "upe_u"
$ws.Range("B14").Value = 0.65
$ws.Range("C14").Value = 0.87
$ws.Range("D14").Value = 0.71
$ws.Range("A15").Value = "upe_l"
$ws.Range("B15").Value = 0.79
$ws.Range("C15").Value = 0.5
$ws.Range("D15").Value = 0.4
$ws.Range("A2").Value = "egp"
$ws.Range("B2").Value = -0.17
$ws.Range("C2").Value = -0.23
$ws.Range("D2").Value = -0.14
$ws.Range("A3").Value = "kan_l"
$ws.Range("B3").Value = 0.48
$ws.Range("C3").Value = 0.51
$ws.Range("D3").Value = 0.69

# formatting: column A + header row centred, data block shown with 2 decimals
$ws.Range("A1:A15").HorizontalAlignment = -4108
$ws.Range("A1:A15").VerticalAlignment = -4108
$ws.Range("B1:D1").HorizontalAlignment = -4108
$ws.Range("B1:D1").VerticalAlignment = -4108
$ws.Range("B2:D15").NumberFormat = "0.00"
$ws.Range("B2:D15").HorizontalAlignment = -4108
$ws.Range("B2:D15").VerticalAlignment = -4108

# red-yellow-green color scale over the data block
$null = $ws.Range("B2:D15").FormatConditions.AddColorScale(3)

# view: zoomed in, this is the sheet/tab that ends up active
$ws.Range("D15").Select()
$ws.Activate()
$excel.ActiveWindow.Zoom = 280
